$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) In the "Bono Anual" row of the remuneration table, prefix the
#    description with "Hasta " so it reads
#    "Hasta 2 S.B. Según Resultados (Contrato Indefinido)".
# ------------------------------------------------------------------
$targetRange = $d.Content
$found = $targetRange.Find.Execute("2 S.B. Según Resultados", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $targetRange.InsertBefore("Hasta ")
}

# ------------------------------------------------------------------
# 2) Word keeps a single "_GoBack" bookmark marking the location of
#    the most recent edit. Move it from its old position (end of the
#    "Estructura Anual de Remuneraciones..." heading) to the end of
#    "Lugar de Trabajo: {{lugar_de_trabajo}}", which is where editing
#    finished. Re-adding a bookmark named "_GoBack" automatically
#    removes the previous one, since Word treats it as a singleton.
# ------------------------------------------------------------------
$bookmarkRange = $d.Content
$bookmarkRange.Find.Execute("Lugar de Trabajo: {{lugar_de_trabajo}}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
